$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23 (shifts old rows 23-33 down to 24-34)
$ws.Rows.Item(23).Insert()

# Fill in the new person: Lauren Macaisa
$ws.Range("A23").Value = "Lauren"
$ws.Range("B23").Value = "Macaisa"
$ws.Range("C23").Value = 5
$ws.Range("F23").Value = "macaisa@broadinstitute.org"
$ws.Range("G23").Value = "Lauren is a Research Associate II focused on organizing and optimizing wet lab protocols for the McCarroll & Macosko BICAN lab efforts. Prior to working at the Broad, Lauren worked at Moffitt Cancer Center developing immunotherapy treatments for breast cancer. She graduated from University of North Florida with her B.S. in Behavioral Neuroscience and is pursuing her M.S. in Biotechnology from Northeastern."
$ws.Range("H23").Value = "Data Generation"

# Renumber the "importance" values (column C) for the rest of the
# "Data Generation" block that followed (now rows 24-28), which shift
# up by one to make room for Lauren at importance 5.
$ws.Range("C24").Value = 6
$ws.Range("C25").Value = 7
$ws.Range("C26").Value = 8
$ws.Range("C27").Value = 9
$ws.Range("C28").Value = 10

# Row 24 (Giovanni Marrero) picked up bold formatting on its first-name cell
$ws.Range("A24").Font.Bold = $true

# Restore the active selection
$ws.Range("A1:H23").Select()

# Refresh the AutoFilter range to cover the newly added row
$ws.Range("A1:H34").AutoFilter()
$ws.Range("A1:H34").AutoFilter()

# Keep the _FilterDatabase defined name in sync with the new filter range
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$H`$34"
  }
}

Write-Output "done"
